$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# ---- 1) Update data rows 2-11 (A-E, G, H) with new content ----
# (F column URLs are set further below together with hyperlink re-creation)

$ws.Range("A2").Value = '2025-11-20 06:26:50'
$ws.Range("B2").Value = '【AI開発】生成AI・RAGシステム構築パートナー募集'
$ws.Range("C2").Value = 'システム開発'
$ws.Range("D2").Value = '1,000 ~ 5,000 円 / 固定'
$ws.Range("E2").Value = '期限情報なし'
$ws.Range("G2").Value = 375
$ws.Range("H2").Value = '🔥AI,Ai ◆開発'

$ws.Range("A3").Value = '2025-11-20 06:26:50'
$ws.Range("B3").Value = '【急募】Web管理システム構築・AI機能実装のプロを探しています'
$ws.Range("C3").Value = 'システム開発'
$ws.Range("D3").Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range("E3").Value = '期限情報なし'
$ws.Range("G3").Value = 360
$ws.Range("H3").Value = '🔥AI,Ai ◇管理'

$ws.Range("A4").Value = '2025-11-20 06:26:50'
$ws.Range("B4").Value = '【急募】案件管理システム開発のフリーランス募集!'
$ws.Range("C4").Value = 'システム開発'
$ws.Range("D4").Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range("E4").Value = '期限情報なし'
$ws.Range("G4").Value = 160
$ws.Range("H4").Value = '◆開発,システム開発 ◇管理'

$ws.Range("A5").Value = '2025-11-20 06:26:50'
$ws.Range("B5").Value = 'webアプリの開発'
$ws.Range("C5").Value = 'システム開発'
$ws.Range("D5").Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range("E5").Value = '期限情報なし'
$ws.Range("G5").Value = 100
$ws.Range("H5").Value = '◆開発 ◇アプリ'

$ws.Range("A6").Value = '2025-11-20 06:26:50'
$ws.Range("B6").Value = '【急募】価格更新サイトにエクセルアップロード後、内容を更新するプログラム作成依頼'
$ws.Range("C6").Value = 'システム開発'
$ws.Range("D6").Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Range("E6").Value = '期限情報なし'
$ws.Range("G6").Value = 38
$ws.Range("H6").Value = '◇サイト'

$ws.Range("A7").Value = '2025-11-20 06:26:50'
$ws.Range("B7").Value = '【急募】WordPressにe-SCOTT決済機能を導入'
$ws.Range("C7").Value = 'システム開発'
$ws.Range("D7").Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Range("E7").Value = '期限情報なし'
$ws.Range("G7").Value = 33
$ws.Range("H7").Value = '○WordPress'

$ws.Range("A8").Value = '2025-11-20 06:26:50'
$ws.Range("B8").Value = '【フルスタックエンジニア】 働きながらスキルアップもできるEC業界で活躍してみませんか?'
$ws.Range("C8").Value = 'システム開発'
$ws.Range("D8").Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range("E8").Value = '期限情報なし'
$ws.Range("G8").Value = 25
$ws.Range("H8").ClearContents()

$ws.Range("A9").Value = '2025-11-20 06:26:50'
$ws.Range("B9").Value = '【急募】Flutterflowの扱えるノーコードエンジニアを探しています!'
$ws.Range("C9").Value = 'システム開発'
$ws.Range("D9").Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range("E9").Value = '期限情報なし'
$ws.Range("G9").Value = 25
$ws.Range("H9").ClearContents()

$ws.Range("A10").Value = '2025-11-20 06:26:50'
$ws.Range("B10").Value = '【SESエンジニア募集】多様なプロジェクトに参画可能!'
$ws.Range("C10").Value = 'システム開発'
$ws.Range("D10").Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range("E10").Value = '期限情報なし'
$ws.Range("G10").Value = 25
$ws.Range("H10").ClearContents()

$ws.Range("A11").Value = '2025-11-20 06:26:50'
$ws.Range("B11").Value = 'Networkエンジニア'
$ws.Range("C11").Value = 'システム開発'
$ws.Range("D11").Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range("E11").Value = '期限情報なし'
$ws.Range("G11").Value = 25
$ws.Range("H11").ClearContents()

# ---- 2) Remove now-obsolete rows 12-15 (shrinks dimension to A1:H11) ----
$ws.Range("A12:H15").EntireRow.Delete()

# ---- 3) Resize columns B, D, H (ColumnWidth offset by the 0.8333... Excel padding constant) ----
$ws.Columns.Item(2).ColumnWidth = 45.166666666666664
$ws.Columns.Item(4).ColumnWidth = 29.166666666666668
$ws.Columns.Item(8).ColumnWidth = 15.166666666666666

# ---- 4) Rebuild hyperlinks on column F for rows 2-11 with the updated target URLs ----
$ws.Hyperlinks.Delete()
$ws.Range("F2").Value = 'https://www.lancers.jp/work/detail/5437447'
$ws.Hyperlinks.Add($ws.Range("F2"), 'https://www.lancers.jp/work/detail/5437447')
$ws.Range("F2").Style = "Hyperlink"

$ws.Range("F3").Value = 'https://www.lancers.jp/work/detail/5437717'
$ws.Hyperlinks.Add($ws.Range("F3"), 'https://www.lancers.jp/work/detail/5437717')
$ws.Range("F3").Style = "Hyperlink"

$ws.Range("F4").Value = 'https://www.lancers.jp/work/detail/5437726'
$ws.Hyperlinks.Add($ws.Range("F4"), 'https://www.lancers.jp/work/detail/5437726')
$ws.Range("F4").Style = "Hyperlink"

$ws.Range("F5").Value = 'https://www.lancers.jp/work/detail/5437832'
$ws.Hyperlinks.Add($ws.Range("F5"), 'https://www.lancers.jp/work/detail/5437832')
$ws.Range("F5").Style = "Hyperlink"

$ws.Range("F6").Value = 'https://www.lancers.jp/work/detail/5437655'
$ws.Hyperlinks.Add($ws.Range("F6"), 'https://www.lancers.jp/work/detail/5437655')
$ws.Range("F6").Style = "Hyperlink"

$ws.Range("F7").Value = 'https://www.lancers.jp/work/detail/5437728'
$ws.Hyperlinks.Add($ws.Range("F7"), 'https://www.lancers.jp/work/detail/5437728')
$ws.Range("F7").Style = "Hyperlink"

$ws.Range("F8").Value = 'https://www.lancers.jp/work/detail/5429335'
$ws.Hyperlinks.Add($ws.Range("F8"), 'https://www.lancers.jp/work/detail/5429335')
$ws.Range("F8").Style = "Hyperlink"

$ws.Range("F9").Value = 'https://www.lancers.jp/work/detail/5437783'
$ws.Hyperlinks.Add($ws.Range("F9"), 'https://www.lancers.jp/work/detail/5437783')
$ws.Range("F9").Style = "Hyperlink"

$ws.Range("F10").Value = 'https://www.lancers.jp/work/detail/5437544'
$ws.Hyperlinks.Add($ws.Range("F10"), 'https://www.lancers.jp/work/detail/5437544')
$ws.Range("F10").Style = "Hyperlink"

$ws.Range("F11").Value = 'https://www.lancers.jp/work/detail/5432661'
$ws.Hyperlinks.Add($ws.Range("F11"), 'https://www.lancers.jp/work/detail/5432661')
$ws.Range("F11").Style = "Hyperlink"

